# feat: add 2022-Q3 data
#
# 1) Insert a brand-new worksheet "2022-Q3" right after "总计" (so the tab
#    order becomes 总计, 2022-Q3, 2022-Q1, 2021-Q3, 2021-Q2) and populate it
#    with the Q3 fund-holding table.
# 2) Insert a new row into "总计" summarising the new quarter, pushing the
#    existing rows down and renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q3" sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Match the page margins used throughout the rest of the workbook.
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Header row — reuse the bordered/bold header style already used on every
# other sheet by copying it across before overwriting the cell text.
$styleDonor = $wb.Worksheets.Item("2022-Q1")
$styleDonor.Range("B1:H1").Copy($q3.Range("B1:H1"))

# Column A's index cells also carry that bordered style on every sheet;
# grab a donor range with enough rows (2021-Q3 has eight) to copy from.
$indexStyleDonor = $wb.Worksheets.Item("2021-Q3")
$indexStyleDonor.Range("A2:A4").Copy($q3.Range("A2:A4"))

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Data rows. B-G are stored as text (even the numeric-looking ones), so we
# force Text format before writing, then clear the applied style back to
# Normal so no stray number-format survives on the cell.
$q3Data = @(
    @{ idx = 0; code = "161724"; name = "招商中证煤炭等权指数（LOF）A"; size = "18.69"; pos = "94.42"; pct = "3.25"; mv = "0.6074"; rank = 10 },
    @{ idx = 1; code = "013596"; name = "招商中证煤炭等权指数（LOF）C"; size = "1.38";  pos = "94.42"; pct = "3.25"; mv = "0.0448"; rank = 10 },
    @{ idx = 2; code = "016347"; name = "招商中证煤炭等权指数（LOF）E"; size = "0.10";  pos = "94.42"; pct = "3.25"; mv = "0.0032"; rank = 10 }
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Range("A$r").Value = $row.idx

    $textRange = $q3.Range("B$r`:G$r")
    $textRange.NumberFormat = "@"
    $q3.Range("B$r").Value = $row.code
    $q3.Range("C$r").Value = $row.name
    $q3.Range("D$r").Value = $row.size
    $q3.Range("E$r").Value = $row.pos
    $q3.Range("F$r").Value = $row.pct
    $q3.Range("G$r").Value = $row.mv
    $textRange.Style = "Normal"

    $q3.Range("H$r").Value = $row.rank
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: update "总计" with the new quarter's summary row
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

# Pick up the bordered/centered style used by the rest of column A by
# copying an already-styled neighbour cell, then overwrite its value.
$total.Range("A3").Copy($total.Range("A2"))
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.66
# The row-insert borrowed formatting for B2:D2 from neighbouring rows;
# these columns carry no explicit style in this sheet, so clear it back.
$total.Range("B2:D2").Style = "Normal"

# Renumber the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# Restore the originally-active sheet so tab selection matches the source.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
